# Automatische test-sync: 2025-07-27 19:30:50
# Adds the new "Testmail #9" log entry (row 11) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover the new row, and
# bumps the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row -----------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A11").Value = "Hoi, hebben jullie al iets gehoord?"
$ws.Range("B11").Value = "mailmind.test@zohomail.eu"
$ws.Range("C11").Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$ws.Range("D11").Value = "Overig"
$ws.Range("E11").Value = "Beste,`nDank voor je bericht. Om je zo goed mogelijk te kunnen helpen, zou je wat meer informatie met ons kunnen delen? Kun je aangeven waar je precies op wacht of waarover je meer wilt weten? Op die manier kunnen we je beter van dienst zijn.`nMet vriendelijke groet,  `n[Naam]  `nE-mailassistent"
$ws.Range("F11").Value = "2025-07-27 19:29:55"
$ws.Range("G11").Value = "Ja"
$ws.Range("H11").Value = "Nee"
$ws.Range("I11").Value = "Ja"
$ws.Range("J11").Value = "Nee"

# --- Extend conditional formatting ranges down to the new row ---------
$ws.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D11"))
$ws.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G11"))
$ws.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H11"))
$ws.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I11"))
$ws.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J11"))

# --- Dashboard sheet: "Overig" count goes from 3 to 4 ------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 4
